# Trade #3 closed at 2026-02-17 19:43:51 - unknown UNKNOWN +0.000%
#
# Record the newly-closed MarketMaking trade (#3) across the workbook:
#   - Summary sheet: roll up capital / P&L / trade counters
#   - Strategy Status sheet: update the MarketMaking strategy row
#   - All Trades & MarketMaking sheets: append the new trade row

$wb = $excel.ActiveWorkbook

$summary   = $wb.Worksheets.Item("Summary")
$status    = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$mmSheet   = $wb.Worksheets.Item("MarketMaking")

# ---------------------------------------------------------------------
# Summary sheet updates
# ---------------------------------------------------------------------
$summary.Range("B3").Value = 1300.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.13      # Total P&L %
$summary.Range("B6").Value = 3         # Total Trades
$summary.Range("B7").Value = 1         # Winning Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet updates (MarketMaking row = row 4)
# ---------------------------------------------------------------------
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 3          # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------
# Append new trade row (#3) to both "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = @(
    3,
    "2026-02-17",
    "19:43:45",
    "MarketMaking",
    "DOWN",
    0.64,
    0.7,
    "CLOSED",
    9.375,
    0.06,
    100.02,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.14
)

foreach ($sheet in @($allTrades, $mmSheet)) {
    $sheet.Range("A4").Value = $newRow[0]
    # B4 holds a literal date-shaped string ("YYYY-MM-DD"); without forcing
    # a Text format first, Excel's COM layer auto-converts it to a date
    # serial number, same as it would in real Excel.
    $sheet.Range("B4").NumberFormat = "@"
    $sheet.Range("B4").Value = $newRow[1]
    $sheet.Range("C4").Value = $newRow[2]
    $sheet.Range("D4").Value = $newRow[3]
    $sheet.Range("E4").Value = $newRow[4]
    $sheet.Range("F4").Value = $newRow[5]
    $sheet.Range("G4").Value = $newRow[6]
    $sheet.Range("H4").Value = $newRow[7]
    $sheet.Range("I4").Value = $newRow[8]
    $sheet.Range("J4").Value = $newRow[9]
    $sheet.Range("K4").Value = $newRow[10]
    $sheet.Range("L4").Value = $newRow[11]
    $sheet.Range("M4").Value = $newRow[12]
    $sheet.Range("N4").Value = $newRow[13]
    $sheet.Range("O4").Value = $newRow[14]
    $sheet.Range("P4").Value = $newRow[15]
    $sheet.Range("Q4").Value = $newRow[16]
}
